$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021312911921985
$ws.Cells.Item(2, 4).Value = 1.026629670511691
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.019741662386892
$ws.Cells.Item(2, 9).Value = 1.031243529606694
$ws.Cells.Item(2, 10).Value = 1.026505085214512
$ws.Cells.Item(2, 11).Value = 1.029451757862291
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.0225840102479
$ws.Cells.Item(2, 14).Value = 1.012926493720188
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022083332400925
$ws.Cells.Item(3, 4).Value = 1.027186393188569
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.0211479559186
$ws.Cells.Item(3, 9).Value = 1.031403391665112
$ws.Cells.Item(3, 10).Value = 1.026914003210218
$ws.Cells.Item(3, 11).Value = 1.029816714781862
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.023794745398038
$ws.Cells.Item(3, 14).Value = 1.013061193460122
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.022581819946383
$ws.Cells.Item(4, 4).Value = 1.027546512791065
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.022057966253385
$ws.Cells.Item(4, 9).Value = 1.031505374949012
$ws.Cells.Item(4, 10).Value = 1.027177909381216
$ws.Cells.Item(4, 11).Value = 1.030052045162372
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.024577707883929
$ws.Cells.Item(4, 14).Value = 1.013148120458946
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022791376909254
$ws.Cells.Item(5, 4).Value = 1.027697877796048
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.022440548350716
$ws.Cells.Item(5, 9).Value = 1.031547899163903
$ws.Cells.Item(5, 10).Value = 1.027288689525006
$ws.Cells.Item(5, 11).Value = 1.030150780762836
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.024906756772645
$ws.Cells.Item(5, 14).Value = 1.013184608651216
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022826561994799
$ws.Cells.Item(6, 4).Value = 1.027723290894079
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.022504786465606
$ws.Cells.Item(6, 9).Value = 1.031555018654708
$ws.Cells.Item(6, 10).Value = 1.027307280247423
$ws.Cells.Item(6, 11).Value = 1.03016734731471
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.024961999238732
$ws.Cells.Item(6, 14).Value = 1.013190731894359
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022584620086298
$ws.Cells.Item(7, 4).Value = 1.027548535454251
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.022063078276159
$ws.Cells.Item(7, 9).Value = 1.03150594453377
$ws.Cells.Item(7, 10).Value = 1.027179390283932
$ws.Cells.Item(7, 11).Value = 1.030053365248415
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.024582105074578
$ws.Cells.Item(7, 14).Value = 1.01314860823593
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021573283935093
$ws.Cells.Item(8, 4).Value = 1.026817841082883
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.020216918563294
$ws.Cells.Item(8, 9).Value = 1.03129785737616
$ws.Cells.Item(8, 10).Value = 1.026643423733058
$ws.Cells.Item(8, 11).Value = 1.029575266363517
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.022993281506576
$ws.Cells.Item(8, 14).Value = 1.012972064198863
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.019791023603591
$ws.Cells.Item(9, 4).Value = 1.025529420051517
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.01696394902566
$ws.Cells.Item(9, 9).Value = 1.030920034080462
$ws.Cells.Item(9, 10).Value = 1.02569371322506
$ws.Cells.Item(9, 11).Value = 1.028726536842608
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.020189901497378
$ws.Cells.Item(9, 14).Value = 1.01265919721807
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018602800169582
$ws.Cells.Item(10, 4).Value = 1.024669976073618
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.014795239047474
$ws.Cells.Item(10, 9).Value = 1.030660684196326
$ws.Cells.Item(10, 10).Value = 1.025057061579061
$ws.Cells.Item(10, 11).Value = 1.028156549983053
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.018318350812416
$ws.Cells.Item(10, 14).Value = 1.012449437880351
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.01808828564498
$ws.Cells.Item(11, 4).Value = 1.02429772419909
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.013856103943291
$ws.Cells.Item(11, 9).Value = 1.030546617035249
$ws.Cells.Item(11, 10).Value = 1.024780557715747
$ws.Cells.Item(11, 11).Value = 1.027908759167091
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.01750728777252
$ws.Cells.Item(11, 14).Value = 1.012358331809703
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.01789717198772
$ws.Cells.Item(12, 4).Value = 1.024159438356607
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.013507252996575
$ws.Cells.Item(12, 9).Value = 1.030503982341987
$ws.Cells.Item(12, 10).Value = 1.024677727689015
$ws.Cells.Item(12, 11).Value = 1.027816571592156
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.017205919025334
$ws.Cells.Item(12, 14).Value = 1.01232444922087
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017938166505007
$ws.Cells.Item(13, 4).Value = 1.02418910179615
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.013582083421655
$ws.Cells.Item(13, 9).Value = 1.030513139617143
$ws.Cells.Item(13, 10).Value = 1.0246997906912
$ws.Cells.Item(13, 11).Value = 1.027836352783686
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.017270568395474
$ws.Cells.Item(13, 14).Value = 1.012331719037316
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.01807248812673
$ws.Cells.Item(14, 4).Value = 1.024286293746007
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.013827268130811
$ws.Cells.Item(14, 9).Value = 1.030543098243302
$ws.Cells.Item(14, 10).Value = 1.024772060284367
$ws.Cells.Item(14, 11).Value = 1.027901141909292
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.017482378674776
$ws.Cells.Item(14, 14).Value = 1.012355531914842
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018155248128849
$ws.Cells.Item(15, 4).Value = 1.024346174970947
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.013978332502378
$ws.Cells.Item(15, 9).Value = 1.030561521628908
$ws.Cells.Item(15, 10).Value = 1.024816571523923
$ws.Cells.Item(15, 11).Value = 1.027941041162832
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.017612868122951
$ws.Cells.Item(15, 14).Value = 1.012370198288797
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018636946750636
$ws.Cells.Item(16, 4).Value = 1.024694679049203
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.014857564554462
$ws.Cells.Item(16, 9).Value = 1.03066821726418
$ws.Cells.Item(16, 10).Value = 1.025075394791108
$ws.Cells.Item(16, 11).Value = 1.028172974390612
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.018372164015686
$ws.Cells.Item(16, 14).Value = 1.012455478428195
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018939102573745
$ws.Cells.Item(17, 4).Value = 1.024913258627048
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.015409062248061
$ws.Cells.Item(17, 9).Value = 1.030734671833354
$ws.Cells.Item(17, 10).Value = 1.025237526022134
$ws.Cells.Item(17, 11).Value = 1.028318197303041
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.018848268785445
$ws.Cells.Item(17, 14).Value = 1.012508897854551
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019115344324647
$ws.Cells.Item(18, 4).Value = 1.025040742000123
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.015730735093065
$ws.Cells.Item(18, 9).Value = 1.030773263147972
$ws.Cells.Item(18, 10).Value = 1.025332014362137
$ws.Cells.Item(18, 11).Value = 1.028402808509501
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.019125908290594
$ws.Cells.Item(18, 14).Value = 1.012540029582198
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019175438074999
$ws.Cells.Item(19, 4).Value = 1.025084208737369
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.015840416249379
$ws.Cells.Item(19, 9).Value = 1.030786392855287
$ws.Cells.Item(19, 10).Value = 1.025364218850494
$ws.Cells.Item(19, 11).Value = 1.028431642635472
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.019220565371273
$ws.Cells.Item(19, 14).Value = 1.012550640125732
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018906684172398
$ws.Cells.Item(20, 4).Value = 1.024889808173989
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.015349892457787
$ws.Cells.Item(20, 9).Value = 1.030727559520905
$ws.Cells.Item(20, 10).Value = 1.025220139153579
$ws.Cells.Item(20, 11).Value = 1.02830262606639
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.018797193928944
$ws.Cells.Item(20, 14).Value = 1.01250316923892
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.018032933751803
$ws.Cells.Item(21, 4).Value = 1.024257673541962
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.013755067772909
$ws.Cells.Item(21, 9).Value = 1.030534283488703
$ws.Cells.Item(21, 10).Value = 1.024750782114917
$ws.Cells.Item(21, 11).Value = 1.027882067180562
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.017420008719095
$ws.Cells.Item(21, 14).Value = 1.012348520766797
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.017483572281225
$ws.Cells.Item(22, 4).Value = 1.023860139638181
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.012752253191733
$ws.Cells.Item(22, 9).Value = 1.030411229751045
$ws.Cells.Item(22, 10).Value = 1.024454960460597
$ws.Cells.Item(22, 11).Value = 1.027616794728863
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.016553515408316
$ws.Cells.Item(22, 14).Value = 1.012251045732159
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017774798903499
$ws.Cells.Item(23, 4).Value = 1.024070887631465
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.013283873548875
$ws.Cells.Item(23, 9).Value = 1.030476608085978
$ws.Cells.Item(23, 10).Value = 1.024611849052743
$ws.Cells.Item(23, 11).Value = 1.027757501078394
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.017012918095746
$ws.Cells.Item(23, 14).Value = 1.012302741920595
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018921332654507
$ws.Cells.Item(24, 4).Value = 1.024900404456538
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.015376628757977
$ws.Cells.Item(24, 9).Value = 1.030730773795617
$ws.Cells.Item(24, 10).Value = 1.025227995778448
$ws.Cells.Item(24, 11).Value = 1.028309662331108
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.018820272656111
$ws.Cells.Item(24, 14).Value = 1.012505757836917
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.020251793457502
$ws.Cells.Item(25, 4).Value = 1.025862600466656
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.017804916300633
$ws.Cells.Item(25, 9).Value = 1.031019028657197
$ws.Cells.Item(25, 10).Value = 1.025939857315227
$ws.Cells.Item(25, 11).Value = 1.028946691203637
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.020915094107655
$ws.Cells.Item(25, 14).Value = 1.012740289871162
